$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "lang_code"
$ws.Range("B1").Value = "code"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "descr"
$ws.Range("E1").Value = "is_active"

# Data rows
$ws.Range("A2").Value = "eng"
$ws.Range("B2").Value = "MNA"
$ws.Range("C2").Value = "Manual Adjudication"
$ws.Range("D2").Value = "Rejection during Manual Adjudication"
$ws.Range("E2").Value = $true

$ws.Range("A3").Value = "eng"
$ws.Range("B3").Value = "CLR"
$ws.Range("C3").Value = "Client Rejection"
$ws.Range("D3").Value = "Rejection in Registration Client"
$ws.Range("E3").Value = $true

$ws.Range("A4").Value = "fra"
$ws.Range("B4").Value = "MNA"
$ws.Range("C4").Value = "Manuel arbitrage"
$ws.Range("D4").Value = "Renvoi en cours de sélection manuelle"
$ws.Range("E4").Value = $true

$ws.Range("A5").Value = "fra"
$ws.Range("B5").Value = "CLR"
$ws.Range("C5").Value = "Rejet de client"
$ws.Range("D5").Value = "Rejet en enregistrement Client"
$ws.Range("E5").Value = $true

# Copy the header's formatting (bold font, borders, centered/top alignment) onto
# column A of each data row, matching the style used in the original file.
$ws.Range("A1").Copy()
$ws.Range("A2:A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false
